# edit.ps1
#
# Replaces the placeholder "Hi" paragraph with a tax-report intake line
# ("1.Name: ... Address: ... EIN/ITIN/SS#: ... 1099/W2 amounts ...") and
# updates the document's default (Normal) paragraph style to use
# Times New Roman 12pt, matching the target OOXML diff.

$d = $word.ActiveDocument

# --- The four text segments that make up the new paragraph -----------------
# (segment 1 gets explicit Times New Roman / 12pt run formatting; the rest
#  inherit formatting from the updated "Normal" style, so no explicit rPr
#  is applied to them)
$seg1 = '1.Name:____________________________________________________________'
$seg2 = '  Address:_________________________________________________________'
$seg3 = '  EIN/ITIN/SS#:_____________________________'
$seg4 = '  1099 Amount Paid: $______________ W2 Amount Paid: $______________'

$tailText = $seg2 + $seg3 + $seg4

# --- Rebuild the paragraph's content ---------------------------------------
# The original paragraph contains a visible run ("Hi") followed by a second,
# text-less run that only carries leftover rPr (Calibri/52). Deleting the
# whole story first (rather than just replacing the "Hi" text) drops that
# orphaned formatting run entirely, so the rebuilt paragraph ends up with
# exactly the runs implied by the target XML.
$d.Content.Delete()
$d.Content.InsertAfter($seg1 + $tailText)

# --- Apply explicit run formatting to segment 1 only ------------------------
$seg1Range = $d.Range(0, $seg1.Length)
$seg1Range.Font.Name = "Times New Roman"
$seg1Range.Font.Size = 12

# --- Update the "Normal" style's default run formatting ---------------------
$normalStyle = $d.Styles("Normal")
$normalStyle.Font.Name = "Times New Roman"
$normalStyle.Font.Size = 12
